$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40

# Columns A and D contain values that look like a date / a plain number
# ("2023-06-11" and "24"). Force them to be stored as text (matching the
# existing rows, which hold Date/Week as literal strings) by briefly
# switching the cell to a text number format, then clearing the format
# back off so no explicit style sticks to the cell (matching row 39 etc.,
# which carry no style attribute).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-11"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "18:05:13"
$ws.Cells.Item($row, 3).Value = "Sunday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "24"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 121339
$ws.Cells.Item($row, 6).Value = 134869
$ws.Cells.Item($row, 7).Value = 160900
$ws.Cells.Item($row, 8).Value = 132011
$ws.Cells.Item($row, 9).Value = 176353
$ws.Cells.Item($row, 10).Value = 114465
$ws.Cells.Item($row, 11).Value = 202100
$ws.Cells.Item($row, 12).Value = 222638
$ws.Cells.Item($row, 13).Value = 173827
$ws.Cells.Item($row, 14).Value = 100031
$ws.Cells.Item($row, 15).Value = 38793
$ws.Cells.Item($row, 16).Value = 34218
$ws.Cells.Item($row, 17).Value = 51221
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 37062
$ws.Cells.Item($row, 20).Value = -1
